$wb = $excel.ActiveWorkbook

$smoke = $wb.Worksheets.Item("smoke")
$appControl = $wb.Worksheets.Item("AppControl")

# Flip the Run Flag column (B18:B25) on the "smoke" sheet from "N" to "Y"
$smoke.Range("B18:B25").Value = "Y"

# Keep "smoke" sheet's own remembered selection as a single cell (A26)
$smoke.Range("A26").Select()

# Add the email address + mailto hyperlink to AppControl!B25 (Email ID)
$appControl.Range("B25").Value = "stiyyagura@enhops.com"
$appControl.Hyperlinks.Add($appControl.Range("B25"), "mailto:stiyyagura@enhops.com")
$appControl.Range("B25").Style = "Hyperlink"
$appControl.Range("B25").WrapText = $true

# Make AppControl the active/selected sheet with A26 selected
$appControl.Activate()
$appControl.Range("A26").Select()
